$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell -> new text value (these cells hold numeric-looking
# text/strings in the original workbook, e.g. "20.1", "16.6", ...).
# Setting .NumberFormat = "@" (Text) first keeps Excel's COM layer from
# auto-coercing the numeric-looking string into a real number, so the
# cell keeps being stored as a shared string, matching the source data.
$updates = [ordered]@{
    "D10" = "20.06"
    "B32" = "16.65"
    "D32" = "19.05"
    "B34" = "21.88"
    "C34" = "44.34"
    "D34" = "66.22"
    "B36" = "87.05"
    "C36" = "12.57"
    "D36" = "99.62"
    "B40" = "13.33"
    "C40" = "41.11"
    "D40" = "54.44"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
